$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Cells.Item(51, 8).Value = 6717  # H51: 6977.7827 -> 6717
$ws.Cells.Item(51, 9).Value = 4599.857  # I51: 5442.4287 -> 4599.857
$ws.Cells.Item(51, 10).Value = 7643.25  # J51: 7649.5 -> 7643.25
$ws.Cells.Item(51, 11).Value = 4599.857  # K51: 5442.4287 -> 4599.857
$ws.Cells.Item(51, 12).Value = 7643.25  # L51: 7649.5 -> 7643.25
$ws.Cells.Item(51, 13).Value = -4115.857  # M51: -4958.4287 -> -4115.857
$ws.Cells.Item(51, 14).Value = -8611.25  # N51: -8617.5 -> -8611.25
# Row 137
$ws.Cells.Item(137, 8).Value = 37040440  # H137: 33336416 -> 37040440
$ws.Cells.Item(137, 9).Value = 83336330  # I137: 62502404 -> 83336330
$ws.Cells.Item(137, 10).Value = 3729  # J137: 3859.7144 -> 3729
$ws.Cells.Item(137, 11).Value = 250008990  # K137: 187507212 -> 250008990
$ws.Cells.Item(137, 12).Value = 11187  # L137: 11579.1432 -> 11187
$ws.Cells.Item(137, 13).Value = -250006440  # M137: -187504662 -> -250006440
$ws.Cells.Item(137, 14).Value = -16287  # N137: -16679.1432 -> -16287

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 11661.889  # H2: 11669.667 -> 11661.889
$ws.Cells.Item(2, 9).Value = 619.625  # I2: 628.375 -> 619.625
$ws.Cells.Item(2, 11).Value = 619.625  # K2: 628.375 -> 619.625
$ws.Cells.Item(2, 13).Value = -506.625  # M2: -515.375 -> -506.625
# Row 34
$ws.Cells.Item(34, 8).Value = 254749.75  # H34: 243799.8 -> 254749.75
$ws.Cells.Item(34, 10).Value = 283333  # J34: 262499.75 -> 283333
$ws.Cells.Item(34, 12).Value = 283333  # L34: 262499.75 -> 283333
$ws.Cells.Item(34, 14).Value = -283875  # N34: -263041.75 -> -283875
# Row 61
$ws.Cells.Item(61, 8).Value = 3821.0444  # H61: 3895.3865 -> 3821.0444
$ws.Cells.Item(61, 9).Value = 3816.9773  # I61: 3892.9534 -> 3816.9773
$ws.Cells.Item(61, 11).Value = 3816.9773  # K61: 3892.9534 -> 3816.9773
$ws.Cells.Item(61, 13).Value = -3604.9773  # M61: -3680.9534 -> -3604.9773
# Row 88
$ws.Cells.Item(88, 8).Value = 2738.25  # H88: 3325.75 -> 2738.25
$ws.Cells.Item(88, 10).Value = 1908.6154  # J88: 2631.6924 -> 1908.6154
$ws.Cells.Item(88, 12).Value = 1908.6154  # L88: 2631.6924 -> 1908.6154
$ws.Cells.Item(88, 14).Value = -2720.6154  # N88: -3443.6924 -> -2720.6154
# Row 91
$ws.Cells.Item(91, 8).Value = 2738.25  # H91: 3325.75 -> 2738.25
$ws.Cells.Item(91, 10).Value = 1908.6154  # J91: 2631.6924 -> 1908.6154
$ws.Cells.Item(91, 12).Value = 1908.6154  # L91: 2631.6924 -> 1908.6154
$ws.Cells.Item(91, 14).Value = -4716.6154  # N91: -5439.6924 -> -4716.6154
# Row 116
$ws.Cells.Item(116, 8).Value = 11661.889  # H116: 11669.667 -> 11661.889
$ws.Cells.Item(116, 9).Value = 619.625  # I116: 628.375 -> 619.625
$ws.Cells.Item(116, 11).Value = 619.625  # K116: 628.375 -> 619.625
$ws.Cells.Item(116, 13).Value = 1674.375  # M116: 1665.625 -> 1674.375
# Row 136
$ws.Cells.Item(136, 8).Value = 3821.0444  # H136: 3895.3865 -> 3821.0444
$ws.Cells.Item(136, 9).Value = 3816.9773  # I136: 3892.9534 -> 3816.9773
$ws.Cells.Item(136, 11).Value = 11450.9319  # K136: 11678.8602 -> 11450.9319
$ws.Cells.Item(136, 13).Value = -8900.9319  # M136: -9128.860199999999 -> -8900.9319

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 11661.889  # H3: 11669.667 -> 11661.889
$ws.Cells.Item(3, 9).Value = 619.625  # I3: 628.375 -> 619.625
$ws.Cells.Item(3, 11).Value = 619.625  # K3: 628.375 -> 619.625
$ws.Cells.Item(3, 13).Value = -505.625  # M3: -514.375 -> -505.625
# Row 53
$ws.Cells.Item(53, 8).Value = 0  # H53: 10000 -> 0
$ws.Cells.Item(53, 10).Value = 0  # J53: 10000 -> 0
$ws.Cells.Item(53, 12).Value = 0  # L53: 10000 -> 0
$ws.Cells.Item(53, 14).Value = $null  # N53: -11148 -> (deleted)
# Row 94
$ws.Cells.Item(94, 8).Value = 708.2632  # H94: 751.6111 -> 708.2632
$ws.Cells.Item(94, 9).Value = 685.7059  # I94: 733.0625 -> 685.7059
$ws.Cells.Item(94, 11).Value = 685.7059  # K94: 733.0625 -> 685.7059
$ws.Cells.Item(94, 13).Value = -234.7059  # M94: -282.0625 -> -234.7059
# Row 134
$ws.Cells.Item(134, 8).Value = 2691.5356  # H134: 2761.2222 -> 2691.5356
$ws.Cells.Item(134, 9).Value = 1306.2084  # I134: 1327.7826 -> 1306.2084
$ws.Cells.Item(134, 11).Value = 3918.6252  # K134: 3983.3478 -> 3918.6252
$ws.Cells.Item(134, 13).Value = -1383.6252  # M134: -1448.3478 -> -1383.6252

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 65408.777  # H31: 81539.92999999999 -> 65408.777
$ws.Cells.Item(31, 9).Value = 8061.5557  # I31: 8361.571 -> 8061.5557
$ws.Cells.Item(31, 10).Value = 122756  # J31: 154718.28 -> 122756
$ws.Cells.Item(31, 11).Value = 8061.5557  # K31: 8361.571 -> 8061.5557
$ws.Cells.Item(31, 12).Value = 122756  # L31: 154718.28 -> 122756
$ws.Cells.Item(31, 13).Value = -7766.5557  # M31: -8066.571 -> -7766.5557
$ws.Cells.Item(31, 14).Value = -123346  # N31: -155308.28 -> -123346
# Row 34
$ws.Cells.Item(34, 8).Value = 65408.777  # H34: 81539.92999999999 -> 65408.777
$ws.Cells.Item(34, 9).Value = 8061.5557  # I34: 8361.571 -> 8061.5557
$ws.Cells.Item(34, 10).Value = 122756  # J34: 154718.28 -> 122756
$ws.Cells.Item(34, 11).Value = 8061.5557  # K34: 8361.571 -> 8061.5557
$ws.Cells.Item(34, 12).Value = 122756  # L34: 154718.28 -> 122756
$ws.Cells.Item(34, 13).Value = -7859.5557  # M34: -8159.571 -> -7859.5557
$ws.Cells.Item(34, 14).Value = -123160  # N34: -155122.28 -> -123160
# Row 58
$ws.Cells.Item(58, 8).Value = 4637.5  # H58: 4448.684 -> 4637.5
$ws.Cells.Item(58, 9).Value = 1842.3846  # I58: 1785.7858 -> 1842.3846
$ws.Cells.Item(58, 11).Value = 1842.3846  # K58: 1785.7858 -> 1842.3846
$ws.Cells.Item(58, 13).Value = -1639.3846  # M58: -1582.7858 -> -1639.3846
# Row 105
$ws.Cells.Item(105, 8).Value = 6278.5713  # H105: 5641.5625 -> 6278.5713
$ws.Cells.Item(105, 9).Value = 6372.1  # I105: 5507.1665 -> 6372.1
$ws.Cells.Item(105, 11).Value = 6372.1  # K105: 5507.1665 -> 6372.1
$ws.Cells.Item(105, 13).Value = -4625.1  # M105: -3760.1665 -> -4625.1
# Row 122
$ws.Cells.Item(122, 8).Value = 8347.177  # H122: 8775.125 -> 8347.177
$ws.Cells.Item(122, 10).Value = 17168.834  # J122: 20302.6 -> 17168.834
$ws.Cells.Item(122, 12).Value = 51506.50199999999  # L122: 60907.8 -> 51506.50199999999
$ws.Cells.Item(122, 14).Value = -56406.50199999999  # N122: -65807.79999999999 -> -56406.50199999999
# Row 136
$ws.Cells.Item(136, 8).Value = 4637.5  # H136: 4448.684 -> 4637.5
$ws.Cells.Item(136, 9).Value = 1842.3846  # I136: 1785.7858 -> 1842.3846
$ws.Cells.Item(136, 11).Value = 5527.1538  # K136: 5357.357400000001 -> 5527.1538
$ws.Cells.Item(136, 13).Value = -2977.1538  # M136: -2807.357400000001 -> -2977.1538

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 3282.4375  # H5: 3441.8823 -> 3282.4375
$ws.Cells.Item(5, 10).Value = 7509  # J5: 7292.4287 -> 7509
$ws.Cells.Item(5, 12).Value = 22527  # L5: 21877.2861 -> 22527
$ws.Cells.Item(5, 14).Value = -22751  # N5: -22101.2861 -> -22751
# Row 68
$ws.Cells.Item(68, 8).Value = 15626343  # H68: 15626405 -> 15626343
$ws.Cells.Item(68, 10).Value = 1249.6666  # J68: 1415 -> 1249.6666
$ws.Cells.Item(68, 12).Value = 3748.9998  # L68: 4245 -> 3748.9998
$ws.Cells.Item(68, 14).Value = -5370.9998  # N68: -5867 -> -5370.9998
# Row 71
$ws.Cells.Item(71, 8).Value = 15626343  # H71: 15626405 -> 15626343
$ws.Cells.Item(71, 10).Value = 1249.6666  # J71: 1415 -> 1249.6666
$ws.Cells.Item(71, 12).Value = 11246.9994  # L71: 12735 -> 11246.9994
$ws.Cells.Item(71, 14).Value = -19358.9994  # N71: -20847 -> -19358.9994
# Row 88
$ws.Cells.Item(88, 8).Value = 16137.223  # H88: 16904.375 -> 16137.223
$ws.Cells.Item(88, 10).Value = 18142.857  # J88: 19500 -> 18142.857
$ws.Cells.Item(88, 12).Value = 54428.571  # L88: 58500 -> 54428.571
$ws.Cells.Item(88, 14).Value = -55284.571  # N88: -59356 -> -55284.571
# Row 91
$ws.Cells.Item(91, 8).Value = 16137.223  # H91: 16904.375 -> 16137.223
$ws.Cells.Item(91, 10).Value = 18142.857  # J91: 19500 -> 18142.857
$ws.Cells.Item(91, 12).Value = 54428.571  # L91: 58500 -> 54428.571
$ws.Cells.Item(91, 14).Value = -57392.571  # N91: -61464 -> -57392.571
# Row 135
$ws.Cells.Item(135, 8).Value = 3282.4375  # H135: 3441.8823 -> 3282.4375
$ws.Cells.Item(135, 10).Value = 7509  # J135: 7292.4287 -> 7509
$ws.Cells.Item(135, 12).Value = 67581  # L135: 65631.85830000001 -> 67581
$ws.Cells.Item(135, 14).Value = -72651  # N135: -70701.85830000001 -> -72651

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Cells.Item(102, 8).Value = 4031.5  # H102: 4318.636 -> 4031.5
$ws.Cells.Item(102, 9).Value = 3063  # I102: 3501 -> 3063
$ws.Cells.Item(102, 11).Value = 3063  # K102: 3501 -> 3063
$ws.Cells.Item(102, 13).Value = -1441  # M102: -1879 -> -1441
# Row 113
$ws.Cells.Item(113, 8).Value = 5734.8184  # H113: 5838.3 -> 5734.8184
$ws.Cells.Item(113, 10).Value = 5574.25  # J113: 5699.143 -> 5574.25
$ws.Cells.Item(113, 12).Value = 5574.25  # L113: 5699.143 -> 5574.25
$ws.Cells.Item(113, 14).Value = -9914.25  # N113: -10039.143 -> -9914.25

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 5129.5386  # H22: 5643 -> 5129.5386
$ws.Cells.Item(22, 9).Value = 1833.3  # I22: 2642.2856 -> 1833.3
$ws.Cells.Item(22, 10).Value = 7189.6875  # J22: 6878.5884 -> 7189.6875
$ws.Cells.Item(22, 11).Value = 1833.3  # K22: 2642.2856 -> 1833.3
$ws.Cells.Item(22, 12).Value = 7189.6875  # L22: 6878.5884 -> 7189.6875
$ws.Cells.Item(22, 13).Value = -1538.3  # M22: -2347.2856 -> -1538.3
$ws.Cells.Item(22, 14).Value = -7779.6875  # N22: -7468.5884 -> -7779.6875
# Row 27
$ws.Cells.Item(27, 8).Value = 5129.5386  # H27: 5643 -> 5129.5386
$ws.Cells.Item(27, 9).Value = 1833.3  # I27: 2642.2856 -> 1833.3
$ws.Cells.Item(27, 10).Value = 7189.6875  # J27: 6878.5884 -> 7189.6875
$ws.Cells.Item(27, 11).Value = 1833.3  # K27: 2642.2856 -> 1833.3
$ws.Cells.Item(27, 12).Value = 7189.6875  # L27: 6878.5884 -> 7189.6875
$ws.Cells.Item(27, 13).Value = -1726.3  # M27: -2535.2856 -> -1726.3
$ws.Cells.Item(27, 14).Value = -7403.6875  # N27: -7092.5884 -> -7403.6875
# Row 40
$ws.Cells.Item(40, 8).Value = 9799.666999999999  # H40: 11000.667 -> 9799.666999999999
$ws.Cells.Item(40, 9).Value = 8415.5  # I40: 8713.857 -> 8415.5
$ws.Cells.Item(40, 10).Value = 15336.333  # J40: 19004.5 -> 15336.333
$ws.Cells.Item(40, 11).Value = 8415.5  # K40: 8713.857 -> 8415.5
$ws.Cells.Item(40, 12).Value = 15336.333  # L40: 19004.5 -> 15336.333
$ws.Cells.Item(40, 13).Value = -8279.5  # M40: -8577.857 -> -8279.5
$ws.Cells.Item(40, 14).Value = -15608.333  # N40: -19276.5 -> -15608.333
# Row 122
$ws.Cells.Item(122, 8).Value = 6267.5557  # H122: 6800.7144 -> 6267.5557
$ws.Cells.Item(122, 9).Value = 4560.6  # I122: 4666.6665 -> 4560.6
$ws.Cells.Item(122, 11).Value = 13681.8  # K122: 13999.9995 -> 13681.8
$ws.Cells.Item(122, 13).Value = -11231.8  # M122: -11549.9995 -> -11231.8
# Row 137
$ws.Cells.Item(137, 8).Value = 67975.2  # H137: 69969.836 -> 67975.2
$ws.Cells.Item(137, 9).Value = 60000  # I137: 0 -> 60000
$ws.Cells.Item(137, 10).Value = 69969  # J137: 69969.836 -> 69969
$ws.Cells.Item(137, 11).Value = 60000  # K137: 0 -> 60000
$ws.Cells.Item(137, 12).Value = 69969  # L137: 69969.836 -> 69969
$ws.Cells.Item(137, 13).Value = -54900  # M137: None -> -54900
$ws.Cells.Item(137, 14).Value = -80169  # N137: -80169.836 -> -80169

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Cells.Item(62, 8).Value = 7900  # H62: 7312.5 -> 7900
$ws.Cells.Item(62, 9).Value = 7900  # I62: 7250 -> 7900
$ws.Cells.Item(62, 10).Value = 0  # J62: 7333.3335 -> 0
$ws.Cells.Item(62, 11).Value = 7900  # K62: 7250 -> 7900
$ws.Cells.Item(62, 12).Value = 0  # L62: 7333.3335 -> 0
$ws.Cells.Item(62, 13).Value = -7276  # M62: -6626 -> -7276
$ws.Cells.Item(62, 14).Value = $null  # N62: -8581.333500000001 -> (deleted)
# Row 65
$ws.Cells.Item(65, 8).Value = 7900  # H65: 7312.5 -> 7900
$ws.Cells.Item(65, 9).Value = 7900  # I65: 7250 -> 7900
$ws.Cells.Item(65, 10).Value = 0  # J65: 7333.3335 -> 0
$ws.Cells.Item(65, 11).Value = 39500  # K65: 36250 -> 39500
$ws.Cells.Item(65, 12).Value = 0  # L65: 36666.6675 -> 0
$ws.Cells.Item(65, 13).Value = -36380  # M65: -33130 -> -36380
$ws.Cells.Item(65, 14).Value = $null  # N65: -42906.6675 -> (deleted)
# Row 100
$ws.Cells.Item(100, 8).Value = 999  # H100: 701 -> 999
$ws.Cells.Item(100, 10).Value = 999  # J100: 701 -> 999
$ws.Cells.Item(100, 12).Value = 1998  # L100: 1402 -> 1998
$ws.Cells.Item(100, 14).Value = -3080  # N100: -2484 -> -3080
# Row 107
$ws.Cells.Item(107, 8).Value = 435.33334  # H107: 458.77777 -> 435.33334
$ws.Cells.Item(107, 10).Value = 420.75  # J107: 462.77777 -> 420.75
$ws.Cells.Item(107, 12).Value = 1262.25  # L107: 1388.33331 -> 1262.25
$ws.Cells.Item(107, 14).Value = -5102.25  # N107: -5228.33331 -> -5102.25
# Row 132
$ws.Cells.Item(132, 8).Value = 6570.7144  # H132: 5840.375 -> 6570.7144
$ws.Cells.Item(132, 9).Value = 4661.8975  # I132: 4225.0454 -> 4661.8975
$ws.Cells.Item(132, 10).Value = 14015.1  # J132: 11763.25 -> 14015.1
$ws.Cells.Item(132, 11).Value = 13985.6925  # K132: 12675.1362 -> 13985.6925
$ws.Cells.Item(132, 12).Value = 42045.3  # L132: 35289.75 -> 42045.3
$ws.Cells.Item(132, 13).Value = -11455.6925  # M132: -10145.1362 -> -11455.6925
$ws.Cells.Item(132, 14).Value = -47105.3  # N132: -40349.75 -> -47105.3
# Row 136
$ws.Cells.Item(136, 8).Value = 2023.5807  # H136: 1973.6666 -> 2023.5807
$ws.Cells.Item(136, 9).Value = 1475.037  # I136: 1456.069 -> 1475.037
$ws.Cells.Item(136, 11).Value = 4425.111  # K136: 4368.207 -> 4425.111
$ws.Cells.Item(136, 13).Value = -1875.111  # M136: -1818.207 -> -1875.111
